# Adds four new localization rows (iOS "Select account" screen) to the
# end of the existing table on sheet "iOS".
#
# Columns: A=ELEMENT_KEY, B=ELEMENT_VALUE, C=FRENCH, D=SPANICH, E=VALUE_TYPE
#
# NOTE: for each row, column B (the value) is written before column A
# (the key) so that the new shared-string entries land in the table in
# the same order as the target workbook (value string first, then key
# string), matching indices 453-460.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 227
$ws.Cells.Item(227, 2).Value = "Select account"
$ws.Cells.Item(227, 1).Value = "SELECT_ACCOUNT_TITLE"
$ws.Cells.Item(227, 5).Value = "device-accessibilityid"

# Row 228
$ws.Cells.Item(228, 2).Value = "Close"
$ws.Cells.Item(228, 1).Value = "SELECT_ACCOUNT_CLOSE_BTN"
$ws.Cells.Item(228, 5).Value = "device-accessibilityid"

# Row 229
$ws.Cells.Item(229, 2).Value = '//XCUIElementTypeOther[@name="CURRENTLY SELECTED"]'
$ws.Cells.Item(229, 1).Value = "CURRENTLY_SELECTED_LABEL"
$ws.Cells.Item(229, 5).Value = "device-xpath"

# Row 230
$ws.Cells.Item(230, 2).Value = '//XCUIElementTypeOther[@name="OTHER ACCOUNTS"]'
$ws.Cells.Item(230, 1).Value = "OTHER_ACCOUNTS_LABEL"
$ws.Cells.Item(230, 5).Value = "device-xpath"

# Update the view selection to match the target workbook's sheet view
# (the table grew from 226 to 230 data rows, so the remembered selection
# moves from B229 to B233).
$ws.Range("B233").Select()
